$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the AutoFilter (and its backing defined name) ---------------
$ws.AutoFilterMode = $false
foreach ($n in $wb.Names) {
    $n.Delete() | Out-Null
}

# --- Drop the last data row (old row 10 no longer exists in the refresh) -
$ws.Rows.Item(10).Delete() | Out-Null

# --- Row 2: header row (weight/car-number headers shift into L:O) -------
$ws.Cells.Item(2,1).Value = 'Initial'
$ws.Cells.Item(2,2).Value = 'Number'
$ws.Cells.Item(2,3).Value = 'Location City'
$ws.Cells.Item(2,4).Value = 'State'
$ws.Cells.Item(2,5).Value = 'Month'
$ws.Cells.Item(2,6).Value = 'Day'
$ws.Cells.Item(2,7).Value = 'Time'
$ws.Cells.Item(2,8).Value = 'Event'
$ws.Cells.Item(2,9).Value = 'Train ID'
$ws.Cells.Item(2,10).Value = 'Destination City'
$ws.Cells.Item(2,11).Value = 'State'
$ws.Cells.Item(2,12).Value = 'Gross Weight'
$ws.Cells.Item(2,13).Value = 'Tare Weight'
$ws.Cells.Item(2,14).Value = 'Net Weight'

# --- Row 3: car not authorized, most fields blank ------------------------
$ws.Cells.Item(3,1).Value = 'BN'
$ws.Cells.Item(3,2).Value = 471547
$ws.Cells.Item(3,3).Value = 'Not authorized to view shipment'
$ws.Cells.Item(3,4).ClearContents() | Out-Null
$ws.Cells.Item(3,5).ClearContents() | Out-Null
$ws.Cells.Item(3,6).ClearContents() | Out-Null
$ws.Cells.Item(3,7).ClearContents() | Out-Null
$ws.Cells.Item(3,8).ClearContents() | Out-Null
$ws.Cells.Item(3,9).ClearContents() | Out-Null
$ws.Cells.Item(3,10).ClearContents() | Out-Null
$ws.Cells.Item(3,11).ClearContents() | Out-Null
$ws.Cells.Item(3,12).Value = 'Not authorized to view shipment'
$ws.Cells.Item(3,13).ClearContents() | Out-Null
$ws.Cells.Item(3,14).ClearContents() | Out-Null

# --- Row 4 ----------------------------------------------------------------
$ws.Cells.Item(4,1).Value = 'CRDX'
$ws.Cells.Item(4,2).Value = 15008
$ws.Cells.Item(4,3).Value = 'DENVER'
$ws.Cells.Item(4,4).Value = 'CO'
$ws.Cells.Item(4,5).Value = 6
$ws.Cells.Item(4,6).Value = 14
$ws.Cells.Item(4,7).Value = 303
$ws.Cells.Item(4,8).Value = 'Arrive In-Transit'
$ws.Cells.Item(4,9).Value = 'HKCKDE'
$ws.Cells.Item(4,10).Value = 'LOVELAND'
$ws.Cells.Item(4,11).Value = 'CO'
$ws.Cells.Item(4,12).Value = 286650
$ws.Cells.Item(4,13).Value = 68700
$ws.Cells.Item(4,14).Value = 217950

# --- Row 5 ----------------------------------------------------------------
$ws.Cells.Item(5,1).Value = 'BNSF'
$ws.Cells.Item(5,2).Value = 468933
$ws.Cells.Item(5,3).Value = 'HOLCOMB'
$ws.Cells.Item(5,4).Value = 'KS'
$ws.Cells.Item(5,5).Value = 6
$ws.Cells.Item(5,6).Value = 13
$ws.Cells.Item(5,7).Value = 702
$ws.Cells.Item(5,8).Value = 'Departure'
$ws.Cells.Item(5,9).Value = 'HKCKDE'
$ws.Cells.Item(5,10).Value = 'LOVELAND'
$ws.Cells.Item(5,11).Value = 'CO'
$ws.Cells.Item(5,12).Value = 234960
$ws.Cells.Item(5,13).Value = 63600
$ws.Cells.Item(5,14).Value = 171360

# --- Row 6 (no Train ID for this event) -----------------------------------
$ws.Cells.Item(6,1).Value = 'CRDX'
$ws.Cells.Item(6,2).Value = 15003
$ws.Cells.Item(6,3).Value = 'JOHNSTOWN'
$ws.Cells.Item(6,4).Value = 'CO'
$ws.Cells.Item(6,5).Value = 6
$ws.Cells.Item(6,6).Value = 12
$ws.Cells.Item(6,7).Value = 1304
$ws.Cells.Item(6,8).Value = 'Placed Actual'
$ws.Cells.Item(6,9).ClearContents() | Out-Null
$ws.Cells.Item(6,10).Value = 'LOVELAND'
$ws.Cells.Item(6,11).Value = 'CO'
$ws.Cells.Item(6,12).Value = 286450
$ws.Cells.Item(6,13).Value = 68400
$ws.Cells.Item(6,14).Value = 218050

# --- Row 7 ----------------------------------------------------------------
$ws.Cells.Item(7,1).Value = 'HRTX'
$ws.Cells.Item(7,2).Value = 541059
$ws.Cells.Item(7,3).Value = 'LITTLETON'
$ws.Cells.Item(7,4).Value = 'CO'
$ws.Cells.Item(7,5).Value = 6
$ws.Cells.Item(7,6).Value = 15
$ws.Cells.Item(7,7).Value = 101
$ws.Cells.Item(7,8).Value = 'Departure'
$ws.Cells.Item(7,9).Value = 'HKCKDE'
$ws.Cells.Item(7,10).Value = 'LOVELAND'
$ws.Cells.Item(7,11).Value = 'CO'
$ws.Cells.Item(7,12).Value = 261250
$ws.Cells.Item(7,13).Value = 64200
$ws.Cells.Item(7,14).Value = 197050

# --- Row 8 ----------------------------------------------------------------
$ws.Cells.Item(8,1).Value = 'CRDX'
$ws.Cells.Item(8,2).Value = 15803
$ws.Cells.Item(8,3).Value = 'LOVELAND'
$ws.Cells.Item(8,4).Value = 'CO'
$ws.Cells.Item(8,5).Value = 6
$ws.Cells.Item(8,6).Value = 12
$ws.Cells.Item(8,7).Value = 1045
$ws.Cells.Item(8,8).Value = 'Junction Received'
$ws.Cells.Item(8,9).Value = 'BNSF'
$ws.Cells.Item(8,10).Value = 'LOVELAND'
$ws.Cells.Item(8,11).Value = 'CO'
$ws.Cells.Item(8,12).Value = 284700
$ws.Cells.Item(8,13).Value = 66900
$ws.Cells.Item(8,14).Value = 217800

# --- Row 9 ----------------------------------------------------------------
$ws.Cells.Item(9,1).Value = 'HRTX'
$ws.Cells.Item(9,2).Value = 541048
$ws.Cells.Item(9,3).Value = 'LOVELAND'
$ws.Cells.Item(9,4).Value = 'CO'
$ws.Cells.Item(9,5).Value = 6
$ws.Cells.Item(9,6).Value = 12
$ws.Cells.Item(9,7).Value = 1045
$ws.Cells.Item(9,8).Value = 'Junction Received'
$ws.Cells.Item(9,9).Value = 'BNSF'
$ws.Cells.Item(9,10).Value = 'LOVELAND'
$ws.Cells.Item(9,11).Value = 'CO'
$ws.Cells.Item(9,12).Value = 202800
$ws.Cells.Item(9,13).Value = 0
$ws.Cells.Item(9,14).Value = 202800

# --- Row 1: refreshed report description/footer -------------------------
$ws.Cells.Item(1,1).Value = 'Description unknown, completed 06/15/2023 05:53:58 EDT, by WPJTOWN1.The search returned: 7 events.'

# --- Column O (Car_no): initial + car number, row by row ------------------
$ws.Cells.Item(3,15).Value = 'BN471547'
$ws.Cells.Item(4,15).Value = 'CRDX15008'
$ws.Cells.Item(5,15).Value = 'BNSF468933'
$ws.Cells.Item(6,15).Value = 'CRDX15003'
$ws.Cells.Item(7,15).Value = 'HRTX541059'
$ws.Cells.Item(8,15).Value = 'CRDX15803'
$ws.Cells.Item(9,15).Value = 'HRTX541048'

# --- Tidy up selection to match the refreshed report ---------------------
$ws.Range("O3:O9").Select() | Out-Null
